$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New DEC codes to add, in order, starting at row 93 (the existing blank
# spacer rows 93-95 get filled in too, continuing the same pattern as the
# DEC_08xx/DEC_09xx rows above them).
$codes = @(
    "DEC_0923","DEC_0924","DEC_0927","DEC_0928","DEC_0929","DEC_0930","DEC_0932","DEC_0933",
    "DEC_0935","DEC_0936","DEC_0937","DEC_0938","DEC_0940","DEC_0941","DEC_0942","DEC_0944",
    "DEC_0945","DEC_0947","DEC_0949","DEC_0950","DEC_0951","DEC_0952","DEC_0954","DEC_0956",
    "DEC_0957","DEC_0959","DEC_0960","DEC_0964","DEC_0965","DEC_0968","DEC_0970","DEC_0971",
    "DEC_0973","DEC_0974","DEC_0977","DEC_0978","DEC_0980","DEC_0981"
)

$startRow = 93
$count = $codes.Count
$lastNewRow = $startRow - 1 + $count   # 130
$spacerRow = $lastNewRow + 1           # 131
$oldRow96NewPos = $spacerRow + 1       # 132 (where the old summary row lands)

# 1) Make room: insert brand-new rows at 96 so that the existing summary block
#    (currently at row 96) ends up at row 132, leaving 93..130 for the new DEC
#    pattern rows and row 131 as the lone blank spacer.
$insertCount = $oldRow96NewPos - 96   # 36
$ws.Rows("96:" + (96 + $insertCount - 1)).Insert()

# 2) Seed rows 93..131 by copying the whole row format+value of row 92 (the
#    last existing DEC pattern row: A=code / B / C / D:J=SIN_DATO) down, one
#    row at a time so every destination row gets its own copy. Row 131 gets
#    the same treatment first and is blanked out afterwards (step 4) so it
#    ends up as the lone spacer row.
$srcRow = $ws.Rows(92)
$srcRow.Copy()
for ($r = $startRow; $r -le $spacerRow; $r++) {
    $ws.Rows($r).PasteSpecial(-4104)
}
$excel.CutCopyMode = 0

# 3) Overwrite column A of each new/changed data row with its DEC code.
for ($i = 0; $i -lt $count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $codes[$i]
}

# 4) Row 131 becomes the lone blank spacer row (same look as the old 93/94/95
#    spacer rows: only B/C styled, no values).
$ws.Range("A" + $spacerRow + ":J" + $spacerRow).ClearContents()

# 5) Update the sheet view to match (scrolled down, selection on the new last
#    data row).
$excel.ActiveWindow.ScrollRow = 112
$ws.Range("B130:J130").Select()
